$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.15000000596046448
$ws.Range("C2").Value = 110
$ws.Range("D2").Value = 0.60000002384185791
$ws.Range("E2").Value = 1.1999998092651367

$ws.Range("B3").Value = 0.20000001788139343
$ws.Range("C3").Value = 80
$ws.Range("D3").Value = 0.70000004768371582
$ws.Range("E3").Value = 1.3499996662139893
